# Fix copy/paste leftovers in the "Short Name" (column B) of the audit-field
# rows across the ERD documentation sheets.
#
# 1) CUST_ADDRESS!B22 had a typo'd short name ("CSTADd_ACPT_TS_UTC_OFST" -
#    lowercase "d") - fix the casing so it reads "CSTADD_ACPT_TS_UTC_OFST".
# 2) CUST_DETAILS!B17:B26 (the Audit Fields block) were still showing the
#    short names copied from the CUST_NAME sheet ("CSTNAME_*") instead of
#    this sheet's own "CSTDET_*" names - correct each one.

$wb = $excel.ActiveWorkbook

$wsAddress = $wb.Worksheets.Item("CUST_ADDRESS")
$wsAddress.Range("B22").Value = "CSTADD_ACPT_TS_UTC_OFST"

$wsDetails = $wb.Worksheets.Item("CUST_DETAILS")
$wsDetails.Range("B17").Value = "CSTDET_EFCTV_DATE"
$wsDetails.Range("B18").Value = "CSTDET_CRUD_VALUE"
$wsDetails.Range("B19").Value = "CSTDET_USER_ID"
$wsDetails.Range("B20").Value = "CSTDET_WS_ID"
$wsDetails.Range("B21").Value = "CSTDET_PRGM_ID"
$wsDetails.Range("B22").Value = "CSTDET_HOST_TS"
$wsDetails.Range("B23").Value = "CSTDET_LOCAL_TS"
$wsDetails.Range("B24").Value = "CSTDET_ACPT_TS"
$wsDetails.Range("B25").Value = "CSTDET_ACPT_TS_UTC_OFST"
$wsDetails.Range("B26").Value = "CSTDET_UUID"
